# Update snowballing with authors' answers
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Mark rows with an asterisk where the authors responded ---
$ws.Range("A3").Value  = "Becker and Chambers**"
$ws.Range("A6").Value  = "CXXR*"
$ws.Range("A8").Value  = "ES3*"
$ws.Range("A9").Value  = "ESSW*"
$ws.Range("A12").Value = "Magni*"
$ws.Range("A13").Value = "Michaelides et al.*"
$ws.Range("A15").Value = "Provenance Curious*"
$ws.Range("A24").Value = "VCR*"

# --- Replace "Tariq, Ali, and Gehani" with a new "SPADE" entry and   ---
# --- shift StarFlow / Sumatra down, picking up the authors' answers ---
$ws.Range("A20").Value = "SPADE"
$ws.Range("A21").Value = "StarFlow"
$ws.Range("A22").Value = "Sumatra*"

$ws.Range("B20").Value = "Comprehension"
$ws.Range("C20").Value = "✗"
$ws.Range("D20").Value = "✓"
$ws.Range("E20").Value = "✗"
$ws.Range("F20").Value = "✗"
$ws.Range("G20").Value = "✓"
$ws.Range("H20").Value = "✗"

$ws.Range("B21").Value = "Management"
$ws.Range("C21").Value = "✓"
$ws.Range("D21").Value = "✓"
$ws.Range("E21").Value = "✗"
$ws.Range("F21").Value = "✓"
$ws.Range("G21").Value = "✓"
$ws.Range("H21").Value = "✗"

$ws.Range("B22").Value = "Reproducibility"
$ws.Range("C22").Value = "✗"
$ws.Range("D22").Value = "✓"
$ws.Range("E22").Value = "✗"
$ws.Range("F22").Value = "✓"
$ws.Range("G22").Value = "✓"
$ws.Range("H22").Value = "✗"

# --- Other individual answer updates ---
$ws.Range("D5").Value  = "✗"
$ws.Range("G7").Value  = "✗"
$ws.Range("B18").Value = "Management"
$ws.Range("G18").Value = "✗"

# --- Recalculated totals ---
$ws.Range("D29").Value = 23
$ws.Range("F30").Value = 7
$ws.Range("G30").Value = 6
